$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-02-11 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-02-12 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("949÷9=105, 4", $true, $false, $false, $false, $false, $true, 1, $false, "339÷7=48, 3", 2) | Out-Null
$d.Content.Find.Execute("877÷7=125, 2", $true, $false, $false, $false, $false, $true, 1, $false, "178÷9=19, 7", 2) | Out-Null
$d.Content.Find.Execute("168÷9=18, 6", $true, $false, $false, $false, $false, $true, 1, $false, "806÷6=134, 2", 2) | Out-Null
$d.Content.Find.Execute("566÷8=70, 6", $true, $false, $false, $false, $false, $true, 1, $false, "669÷5=133, 4", 2) | Out-Null
$d.Content.Find.Execute("577÷2=288, 1", $true, $false, $false, $false, $false, $true, 1, $false, "884÷2=442, 0", 2) | Out-Null
$d.Content.Find.Execute("824÷8=103, 0", $true, $false, $false, $false, $false, $true, 1, $false, "540÷5=108, 0", 2) | Out-Null
$d.Content.Find.Execute("346÷6=57, 4", $true, $false, $false, $false, $false, $true, 1, $false, "251÷4=62, 3", 2) | Out-Null
$d.Content.Find.Execute("520÷4=130, 0", $true, $false, $false, $false, $false, $true, 1, $false, "777÷5=155, 2", 2) | Out-Null
$d.Content.Find.Execute("443÷9=49, 2", $true, $false, $false, $false, $false, $true, 1, $false, "675÷6=112, 3", 2) | Out-Null
$d.Content.Find.Execute("791÷9=87, 8", $true, $false, $false, $false, $false, $true, 1, $false, "950÷6=158, 2", 2) | Out-Null
$d.Content.Find.Execute("545÷4=136, 1", $true, $false, $false, $false, $false, $true, 1, $false, "634÷5=126, 4", 2) | Out-Null
$d.Content.Find.Execute("610÷7=87, 1", $true, $false, $false, $false, $false, $true, 1, $false, "293÷6=48, 5", 2) | Out-Null
$d.Content.Find.Execute("482÷7=68, 6", $true, $false, $false, $false, $false, $true, 1, $false, "505÷5=101, 0", 2) | Out-Null
$d.Content.Find.Execute("712÷8=89, 0", $true, $false, $false, $false, $false, $true, 1, $false, "345÷5=69, 0", 2) | Out-Null
$d.Content.Find.Execute("152÷4=38, 0", $true, $false, $false, $false, $false, $true, 1, $false, "721÷3=240, 1", 2) | Out-Null
$d.Content.Find.Execute("247÷7=35, 2", $true, $false, $false, $false, $false, $true, 1, $false, "749÷5=149, 4", 2) | Out-Null
$d.Content.Find.Execute("862÷5=172, 2", $true, $false, $false, $false, $false, $true, 1, $false, "921÷8=115, 1", 2) | Out-Null
$d.Content.Find.Execute("363÷5=72, 3", $true, $false, $false, $false, $false, $true, 1, $false, "886÷8=110, 6", 2) | Out-Null
$d.Content.Find.Execute("949÷2=474, 1", $true, $false, $false, $false, $false, $true, 1, $false, "496÷9=55, 1", 2) | Out-Null
$d.Content.Find.Execute("809÷8=101, 1", $true, $false, $false, $false, $false, $true, 1, $false, "196÷9=21, 7", 2) | Out-Null
$d.Content.Find.Execute("576÷2=288, 0", $true, $false, $false, $false, $false, $true, 1, $false, "305÷8=38, 1", 2) | Out-Null
$d.Content.Find.Execute("172÷4=43, 0", $true, $false, $false, $false, $false, $true, 1, $false, "252÷7=36, 0", 2) | Out-Null
$d.Content.Find.Execute("885÷7=126, 3", $true, $false, $false, $false, $false, $true, 1, $false, "106÷5=21, 1", 2) | Out-Null
$d.Content.Find.Execute("121÷5=24, 1", $true, $false, $false, $false, $false, $true, 1, $false, "523÷6=87, 1", 2) | Out-Null
$d.Content.Find.Execute("984÷5=196, 4", $true, $false, $false, $false, $false, $true, 1, $false, "199÷6=33, 1", 2) | Out-Null
